$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header relabeling ---
# B2 used to read "unnamed: 1_level_1" and F2 used to read "unnamed: 5_level_1";
# both are corrected to "total". C2/D2/E2 already show the correct text
# ("total", "condição de ocupação na semana de referência",
# "condição de ocupação na semana de referência.1") and are left untouched.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# --- Remove the label-only separator rows that had no data next to them
# ("situação do domicílio" in row 5 and "grandes regiões" in row 8). Deleting
# them shifts the data rows below up so each data row lines up with its
# correct area/region label (urbana, rural, norte, nordeste, sudeste, sul). ---
$ws.Rows("5").Delete()
# After the first deletion, the former row 8 ("grandes regiões") is now row 7.
$ws.Rows("7").Delete()

# Resulting layout (rows 1:10, matching the new dimension A1:F10):
#  1: column headers (unchanged)
#  2: relabeled sub-headers
#  3: (blank, as before)
#  4: brasil
#  5: urbana
#  6: rural
#  7: norte
#  8: nordeste
#  9: sudeste
# 10: sul
